$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the auto date placeholder text ("5/7/21" -> "5/14/21") on the
#    slide master, every slide layout, and the notes master.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq "5/7/21") {
                    $shp.TextFrame.TextRange.Text = "5/14/21"
                }
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout attached to the master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Notes master
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes

# ---------------------------------------------------------------------------
# 2) Slide 3 ("DRY: Don't Repeat Yourself"): merge the trailing "Individual "
#    paragraph into the previous paragraph, lower-casing it to "individual "
#    so the sentence reads "...In this case DRY applies within the
#    individual modules."
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$contentShape = $slide3.Shapes.Item(2)
$tf = $contentShape.TextFrame
$tr = $tf.TextRange

$paragraphs = $tr.Paragraphs()
for ($pi = 1; $pi -le $paragraphs.Count; $pi++) {
    $para = $paragraphs.Item($pi)
    if ($para.Text -eq "Individual ") {
        $prev = $paragraphs.Item($pi - 1)
        if ($prev.Text -like "*In this case DRY applies within the *") {
            $prev.Text = $prev.Text -replace "within the $", "within the individual "
            $para.Delete()
        }
    }
}
